# Updated cryptos list on Thu Jun 27 13:36:27 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table, and
# swaps the ranking of "Stacks" and "FirstDigitalUSD" (rows 43 & 44).
#
# All of these source cells are stored as literal TEXT (not numbers/percentages)
# in the workbook, even when the text happens to look like a plain decimal
# (e.g. "0.997", "7.69"). A naive `Range.Value = "0.997"` assignment lets
# Excel's type-inference turn that into a real number (and, for some of the
# percentage-looking "price" strings, even re-render it in scientific
# notation), which would silently change both the stored cell type and the
# number format/style applied to the cell. To avoid that, every write goes
# through Set-TextValue below, which stuffs the value into the cell via a
# `="literal"` formula (guaranteeing Excel treats it as text) and then
# immediately collapses that formula down to its literal text result with a
# copy / paste-values round-trip - leaving a plain text cell with no leftover
# formula and no style change, exactly like the original inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $formula = '="' + $val + '"'
    $cell.Formula = $formula
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false

Set-TextValue $ws.Cells.Item(2,4) "61.753.43"
Set-TextValue $ws.Cells.Item(2,5) "  +0.43%  "
Set-TextValue $ws.Cells.Item(3,4) "3.457.71"
Set-TextValue $ws.Cells.Item(3,5) "  +2.44%  "
Set-TextValue $ws.Cells.Item(4,4) "0.997"
Set-TextValue $ws.Cells.Item(4,5) "  -0.18%  "
Set-TextValue $ws.Cells.Item(5,4) "576.83"
Set-TextValue $ws.Cells.Item(5,5) "  +0.69%  "
Set-TextValue $ws.Cells.Item(6,4) "147.90"
Set-TextValue $ws.Cells.Item(6,5) "  +8.07%  "
Set-TextValue $ws.Cells.Item(7,4) "3.454.20"
Set-TextValue $ws.Cells.Item(7,5) "  +2.40%  "
Set-TextValue $ws.Cells.Item(8,5) "  +0.01%  "
Set-TextValue $ws.Cells.Item(9,5) "  +0.14%  "
Set-TextValue $ws.Cells.Item(10,4) "7.69"
Set-TextValue $ws.Cells.Item(10,5) "  +3.43%  "
Set-TextValue $ws.Cells.Item(11,5) "  +0.06%  "
Set-TextValue $ws.Cells.Item(12,4) "0.392"
Set-TextValue $ws.Cells.Item(12,5) "  +0.24%  "
Set-TextValue $ws.Cells.Item(13,4) "4.031.18"
Set-TextValue $ws.Cells.Item(13,5) "  +2.09%  "
Set-TextValue $ws.Cells.Item(14,5) "  -1.23%  "
Set-TextValue $ws.Cells.Item(15,4) "27.17"
Set-TextValue $ws.Cells.Item(15,5) "  +4.38%  "
Set-TextValue $ws.Cells.Item(16,5) "  +0.16%  "
Set-TextValue $ws.Cells.Item(17,4) "3.442.70"
Set-TextValue $ws.Cells.Item(17,5) "  +2.06%  "
Set-TextValue $ws.Cells.Item(18,4) "61.738.55"
Set-TextValue $ws.Cells.Item(18,5) "  +0.28%  "
Set-TextValue $ws.Cells.Item(19,4) "6.15"
Set-TextValue $ws.Cells.Item(19,5) "  +4.29%  "
Set-TextValue $ws.Cells.Item(20,4) "14.07"
Set-TextValue $ws.Cells.Item(20,5) "  +0.78%  "
Set-TextValue $ws.Cells.Item(21,5) "  +2.08%  "
Set-TextValue $ws.Cells.Item(22,4) "383.42"
Set-TextValue $ws.Cells.Item(22,5) "  +1.54%  "
Set-TextValue $ws.Cells.Item(23,4) "0.562"
Set-TextValue $ws.Cells.Item(23,5) "  +1.27%  "
Set-TextValue $ws.Cells.Item(24,4) "3.571.78"
Set-TextValue $ws.Cells.Item(24,5) "  +1.72%  "
Set-TextValue $ws.Cells.Item(25,4) "0.995"
Set-TextValue $ws.Cells.Item(25,5) "  -0.55%  "
Set-TextValue $ws.Cells.Item(26,4) "71.97"
Set-TextValue $ws.Cells.Item(26,5) "  +0.78%  "
Set-TextValue $ws.Cells.Item(27,4) "0.0000125"
Set-TextValue $ws.Cells.Item(27,5) "  -0.13%  "
Set-TextValue $ws.Cells.Item(28,4) "0.177"
Set-TextValue $ws.Cells.Item(28,5) "  +9.71%  "
Set-TextValue $ws.Cells.Item(29,4) "7.75"
Set-TextValue $ws.Cells.Item(29,5) "  +2.79%  "
Set-TextValue $ws.Cells.Item(30,4) "1.60"
Set-TextValue $ws.Cells.Item(30,5) "  -8.92%  "
Set-TextValue $ws.Cells.Item(31,5) "  -0.09%  "
Set-TextValue $ws.Cells.Item(32,4) "8.22"
Set-TextValue $ws.Cells.Item(32,5) "  -0.34%  "
Set-TextValue $ws.Cells.Item(33,5) "  +0.44%  "
Set-TextValue $ws.Cells.Item(34,5) "  -0.11%  "
Set-TextValue $ws.Cells.Item(35,4) "23.97"
Set-TextValue $ws.Cells.Item(35,5) "  +1.70%  "
Set-TextValue $ws.Cells.Item(36,4) "5.28"
Set-TextValue $ws.Cells.Item(36,5) "  +1.06%  "
Set-TextValue $ws.Cells.Item(37,4) "7.01"
Set-TextValue $ws.Cells.Item(37,5) "  +2.62%  "
Set-TextValue $ws.Cells.Item(38,5) "  +2.64%  "
Set-TextValue $ws.Cells.Item(39,4) "166.68"
Set-TextValue $ws.Cells.Item(39,5) "  +1.15%  "
Set-TextValue $ws.Cells.Item(40,4) "0.0792"
Set-TextValue $ws.Cells.Item(40,5) "  +2.78%  "
Set-TextValue $ws.Cells.Item(41,4) "26.31"
Set-TextValue $ws.Cells.Item(41,5) "  +7.05%  "
Set-TextValue $ws.Cells.Item(42,4) "0.792"
Set-TextValue $ws.Cells.Item(42,5) "  +2.54%  "

# Rows 43 & 44 swap places in the ranking: "Stacks" drops below
# "FirstDigitalUSD". The rank index in column A stays put; B-E content swaps.
Set-TextValue $ws.Cells.Item(43,2) "FirstDigitalUSD"
Set-TextValue $ws.Cells.Item(43,3) "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Cells.Item(43,4) "1.00"
Set-TextValue $ws.Cells.Item(43,5) "  +0.24%  "
Set-TextValue $ws.Cells.Item(44,2) "Stacks"
Set-TextValue $ws.Cells.Item(44,3) "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Cells.Item(44,4) "1.73"
Set-TextValue $ws.Cells.Item(44,5) "  +0.77%  "

Set-TextValue $ws.Cells.Item(45,4) "42.24"
Set-TextValue $ws.Cells.Item(45,5) "  +1.52%  "
Set-TextValue $ws.Cells.Item(46,5) "  +1.97%  "
Set-TextValue $ws.Cells.Item(47,5) "  -0.30%  "
Set-TextValue $ws.Cells.Item(48,4) "2.659.72"
Set-TextValue $ws.Cells.Item(48,5) "  +13.08%  "
Set-TextValue $ws.Cells.Item(49,4) "23.91"
Set-TextValue $ws.Cells.Item(49,5) "  +4.84%  "
Set-TextValue $ws.Cells.Item(50,5) "  +0.27%  "
Set-TextValue $ws.Cells.Item(51,4) "2.19"
Set-TextValue $ws.Cells.Item(51,5) "  +7.77%  "

$excel.CutCopyMode = $false
